{"js": "// Update the CLIN* custom styles used by this clinical-context report\n// template: switch the theme font from Calibri to Aptos and tweak the\n// point sizes for the heading / sub-heading / bullet / body styles\n// (MHD2-259: report template changes for reporting on 136 genes).\n\nconst styles = context.document.getStyles();\n\n// style display name (\"nameLocal\") -> new font size in points\n// (OOXML half-points / 2). Each entry also records the underlying\n// styleId purely as a fallback lookup key.\nconst updates = [\n  { id: \"CLIN1HEADING\", name: \"CLIN1(HEADING)\", size: 15 },             // <w:sz w:val=\"30\"/> (newly added)\n  { id: \"CLIN2SUBHEADINGS\", name: \"CLIN2(SUBHEADINGS)\", size: 10 },     // 22 -> 20 half-points\n  { id: \"CLIN1HEADINGChar\", name: \"CLIN1(HEADING) Char\", size: 15 },    // 32 -> 30 half-points\n  { id: \"CLIN3BULLETPOINTS\", name: \"CLIN3(BULLET POINTS)\", size: 8 },   // <w:sz w:val=\"16\"/> (newly added)\n  { id: \"CLIN2SUBHEADINGSChar\", name: \"CLIN2(SUBHEADINGS) Char\", size: 10 }, // 26 -> 20 half-points\n  { id: \"CLIN4\", name: \"CLIN4\", size: 5.5 },                            // 12 -> 11 half-points\n  { id: \"CLIN3BULLETPOINTSChar\", name: \"CLIN3(BULLET POINTS) Char\", size: 8 }, // 18 -> 16 half-points\n  { id: \"CLIN4Char\", name: \"CLIN4 Char\", size: 5.5 },                   // 12 -> 11 half-points\n];\n\nconst byName = updates.map((u) => styles.getByNameOrNullObject(u.name));\nconst byId = updates.map((u) => styles.getByNameOrNullObject(u.id));\nfor (const s of byName) s.load(\"isNullObject\");\nfor (const s of byId) s.load(\"isNullObject\");\nawait context.sync();\n\nfor (let i = 0; i < updates.length; i++) {\n  // Prefer the lookup by display (\"local\") name, per the documented\n  // Word.StyleCollection.getByNameOrNullObject behaviour; fall back to\n  // the internal styleId if that one didn't resolve.\n  const style = !byName[i].isNullObject ? byName[i] : byId[i];\n  if (style.isNullObject) continue;\n  const font = style.font;\n  font.name = \"Aptos\";\n  font.size = updates[i].size;\n}\nawait context.sync();\n", "ps1": "# Update the CLIN* custom styles used by this clinical-context report\n# template: switch the theme font from Calibri to Aptos and tweak the\n# point sizes for the heading / sub-heading / bullet / body styles\n# (MHD2-259: report template changes for reporting on 136 genes).\n\n$d = $word.ActiveDocument\n\n# Each row: internal styleId, display (\"local\") name, new font size in\n# points (OOXML half-points / 2). Styles are looked up by local name\n# (the normal Word.Styles(name) convention) with the styleId kept as a\n# fallback in case the display name can't be resolved.\n$updates = @(\n    @{ Id = \"CLIN1HEADING\";          Name = \"CLIN1(HEADING)\";           Size = 15  }  # <w:sz w:val=\"30\"/> (newly added)\n    @{ Id = \"CLIN2SUBHEADINGS\";      Name = \"CLIN2(SUBHEADINGS)\";       Size = 10  }  # 22 -> 20 half-points\n    @{ Id = \"CLIN1HEADINGChar\";      Name = \"CLIN1(HEADING) Char\";      Size = 15  }  # 32 -> 30 half-points\n    @{ Id = \"CLIN3BULLETPOINTS\";     Name = \"CLIN3(BULLET POINTS)\";     Size = 8   }  # <w:sz w:val=\"16\"/> (newly added)\n    @{ Id = \"CLIN2SUBHEADINGSChar\";  Name = \"CLIN2(SUBHEADINGS) Char\";  Size = 10  }  # 26 -> 20 half-points\n    @{ Id = \"CLIN4\";                 Name = \"CLIN4\";                    Size = 5.5 }  # 12 -> 11 half-points\n    @{ Id = \"CLIN3BULLETPOINTSChar\"; Name = \"CLIN3(BULLET POINTS) Char\"; Size = 8  }  # 18 -> 16 half-points\n    @{ Id = \"CLIN4Char\";             Name = \"CLIN4 Char\";               Size = 5.5 }  # 12 -> 11 half-points\n)\n\nforeach ($u in $updates) {\n    $style = $null\n    try { $style = $d.Styles($u.Name) } catch { $style = $null }\n    if ($null -eq $style) { $style = $d.Styles($u.Id) }\n\n    $style.Font.Name = \"Aptos\"\n    $style.Font.Size = $u.Size\n}\n"}
